# Weekly update: insert a new weekly price record for
# "Feria Lagunitas de Puerto Montt - Espárragos" as the new row 10,
# pushing the previously-existing rows 10-36 down to 11-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 10 (old rows 10..36 shift to 11..37).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Cells.Item(10, 1).Value  = 4
$ws.Cells.Item(10, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value  = "Los Lagos"
$ws.Cells.Item(10, 4).Value  = 44525
$ws.Cells.Item(10, 5).Value  = 10
$ws.Cells.Item(10, 6).Value  = 300000000
$ws.Cells.Item(10, 7).Value  = "Espárragos"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 180
$ws.Cells.Item(10, 11).Value = 1600
$ws.Cells.Item(10, 12).Value = 1600
$ws.Cells.Item(10, 13).Value = 1600
$ws.Cells.Item(10, 14).Value = "$/kilo"
$ws.Cells.Item(10, 15).Value = "Provincia de Linares"
$ws.Cells.Item(10, 16).Value = 1600
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
